$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 219, shifting existing rows 219:236 down to 220:237
$ws.Rows("219:219").Insert()

# Fill the new row 219 with data
$ws.Range("A219").Value = 4
$ws.Range("B219").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C219").Value = "Los Lagos"
$ws.Range("D219").Value = 44578
$ws.Range("E219").Value = 10
$ws.Range("F219").Value = 100114014
$ws.Range("G219").Value = "Betarraga"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 500
$ws.Range("K219").Value = 1000
$ws.Range("L219").Value = 1000
$ws.Range("M219").Value = 1000
$ws.Range("N219").Value = "$/paquete 5 unidades"
$ws.Range("O219").Value = "Región del Maule"
$ws.Range("P219").Value = 200
$ws.Range("Q219").Value = 5
$ws.Range("R219").Value = "Hortaliza"
